$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1779566666666667
$ws.Range("H2").Value = 0.5338700000000001
$ws.Range("I2").Value = 0.01192558037548992
$ws.Range("J2").Value = 0.01192558037548992
$ws.Range("M2").Value = 3.626135
$ws.Range("N2").Value = 10.878405
$ws.Range("O2").Value = 0.4728835835086186
$ws.Range("P2").Value = 0.4728835835086186
$ws.Range("Q2").Value = 0.6452948974833334
$ws.Range("R2").Value = 5.807654077350001
$ws.Range("S2").Value = 0.005639411183381732
$ws.Range("T2").Value = 0.005639411183381732

$ws.Range("G3").Value = 0.1779566666666667
$ws.Range("H3").Value = 0.5338700000000001
$ws.Range("I3").Value = 0.01192558037548992
$ws.Range("J3").Value = 0.01192558037548992
$ws.Range("O3").Value = 0.01581792773244636
$ws.Range("P3").Value = 0.01581792773244636
$ws.Range("Q3").Value = 0.02158507592666667
$ws.Range("R3").Value = 0.19426568334
$ws.Range("S3").Value = 0.0001886379685469801
$ws.Range("T3").Value = 0.0001886379685469801

$ws.Range("G4").Value = 0.1779566666666667
$ws.Range("H4").Value = 0.5338700000000001
$ws.Range("I4").Value = 0.01192558037548992
$ws.Range("J4").Value = 0.01192558037548992
$ws.Range("O4").Value = 0.511298488758935
$ws.Range("P4").Value = 0.511298488758935
$ws.Range("Q4").Value = 0.6977157114211111
$ws.Range("R4").Value = 6.279441402790001
$ws.Range("S4").Value = 0.006097531223561211
$ws.Range("T4").Value = 0.00609753122356121

$ws.Range("I5").Value = 0.540575811616083
$ws.Range("J5").Value = 0.540575811616083
$ws.Range("M5").Value = 3.626135
$ws.Range("N5").Value = 10.878405
$ws.Range("O5").Value = 0.4728835835086186
$ws.Range("P5").Value = 0.4728835835086186
$ws.Range("Q5").Value = 29.25063619173667
$ws.Range("R5").Value = 263.25572572563
$ws.Range("S5").Value = 0.2556294269550933
$ws.Range("T5").Value = 0.2556294269550933

$ws.Range("I6").Value = 0.540575811616083
$ws.Range("J6").Value = 0.540575811616083
$ws.Range("O6").Value = 0.01581792773244636
$ws.Range("P6").Value = 0.01581792773244636
$ws.Range("S6").Value = 0.008550789122051738
$ws.Range("T6").Value = 0.008550789122051738

$ws.Range("I7").Value = 0.540575811616083
$ws.Range("J7").Value = 0.540575811616083
$ws.Range("O7").Value = 0.511298488758935
$ws.Range("P7").Value = 0.511298488758935
$ws.Range("S7").Value = 0.276395595538938
$ws.Range("T7").Value = 0.276395595538938

$ws.Range("G8").Value = 6.677692666666666
$ws.Range("I8").Value = 0.4474986080084269
$ws.Range("J8").Value = 0.4474986080084269
$ws.Range("M8").Value = 3.626135
$ws.Range("N8").Value = 10.878405
$ws.Range("O8").Value = 0.4728835835086186
$ws.Range("P8").Value = 0.4728835835086186
$ws.Range("Q8").Value = 24.21421509784333
$ws.Range("R8").Value = 217.92793588059
$ws.Range("S8").Value = 0.2116147453701435
$ws.Range("T8").Value = 0.2116147453701435

$ws.Range("G9").Value = 6.677692666666666
$ws.Range("I9").Value = 0.4474986080084269
$ws.Range("J9").Value = 0.4474986080084269
$ws.Range("O9").Value = 0.01581792773244636
$ws.Range("P9").Value = 0.01581792773244636
$ws.Range("Q9").Value = 0.8099640543106666
$ws.Range("R9").Value = 7.289676488795999
$ws.Range("S9").Value = 0.007078500641847639
$ws.Range("T9").Value = 0.007078500641847638

$ws.Range("G10").Value = 6.677692666666666
$ws.Range("I10").Value = 0.4474986080084269
$ws.Range("J10").Value = 0.4474986080084269
$ws.Range("O10").Value = 0.511298488758935
$ws.Range("P10").Value = 0.511298488758935
$ws.Range("Q10").Value = 26.18126747845844
$ws.Range("R10").Value = 235.631407306126
$ws.Range("S10").Value = 0.2288053619964358
$ws.Range("T10").Value = 0.2288053619964357

Write-Output "Applied NATMI TPM updates"